# "Basic level design implemented"
#
# 1) Remove the duplicate "Implement Camera" list item that used to sit in
#    the "To be done" section, right before "Implement second player".
# 2) Add a new "Implement Camera" list item to the "Done" section, right
#    after "Implement Player" (same numbering list as the rest of "Done",
#    numId 5), and restructure the trailing placeholder/bookmark paragraph
#    accordingly.

$d = $word.ActiveDocument

# --- Step 1: delete the stray "Implement Camera" paragraph in "To be done" ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Implement Camera`r") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: locate "Implement Player" in the "Done" section ---
$implPlayerIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Implement Player`r") {
        $implPlayerIdx = $i
        break
    }
}

if ($implPlayerIdx -ne -1) {
    $pImplPlayer = $d.Paragraphs.Item($implPlayerIdx)

    # Insert a brand-new paragraph right after it; the runtime clones the
    # "Implement Player" paragraph formatting (ListParagraph style, numId 5
    # list numbering, en-US language) onto the new paragraph automatically.
    $pImplPlayer.Range.InsertParagraphAfter() | Out-Null
    $pNewCamera = $d.Paragraphs.Item($implPlayerIdx + 1)
    $pNewCamera.Range.Text = "Implement Camera"

    # The paragraph that follows is the old placeholder that carries the
    # "_GoBack" bookmark and currently has <w:ind w:left="360"/>.
    $pBookmark = $d.Paragraphs.Item($implPlayerIdx + 2)

    # Insert a fresh empty paragraph after the bookmark paragraph while it
    # still has its original "ind left=360" formatting, so the new
    # paragraph inherits that exact formatting (matches the paragraph that
    # used to sit there before the bookmark paragraph's own style changes).
    $rngAfterBookmark = $d.Range($pBookmark.Range.End, $pBookmark.Range.End)
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rngAfterBookmark.InsertXML($xml) | Out-Null

    # Now convert the bookmark paragraph itself to use the ListParagraph
    # style (dropping its "ind left=360") while keeping its bookmark intact
    # and restoring the en-US language mark that the style switch clears.
    $pBookmark.Style = "List Paragraph"
    $pBookmark.Range.LanguageID = "en-US"
}
